$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (the surviving order record) ---
$ws.Range("A2").Value = 675

# B2 held a date (with a date-number-format style). The new value is a
# plain text date string, so clear the inherited date formatting first,
# then write the text.
$ws.Range("B2").ClearFormats()
$ws.Range("B2").Value = "20-04-25"

# C2 (Cashier) is unchanged.
# D2 (KOT) changes from 8 to 1.
$ws.Range("D2").Value = 1

# E2:H2 (Price / SGST / CGST / Tax) are unchanged.

# I2 (Food Items) changes.
$ws.Range("I2").Value = "Chicken Wrap (x1)"

# --- Remove the old rows 3-9 (only one order record remains) ---
$ws.Range("A3:I9").Delete()
